{"js": "// Strike through the \"casual\" (Member/Visitor) use-case lines that are no\n// longer applicable, and tidy up the \"Facebook\" spell-checked word (which\n// used to be split into two runs with the _GoBack bookmark wedged between\n// them) into a single run, moving that bookmark onto \"Maken profiel\" instead.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Paragraphs (0-based) whose entire text gets struck through.\nconst fullyStruckThrough = [\n  \"Betalen contributie,  Fully-Dressed\",\n  \"Huren golfveld, \",\n  \"Opvragen statistieken,  Fully-Dressed \",\n  \"Opvragen Evenementen\",\n  \"Kopen artikelen,  Fully-Dressed \",\n];\n\nfor (const para of paragraphs.items) {\n  if (fullyStruckThrough.includes(para.text)) {\n    para.font.strikeThrough = true;\n  }\n}\nawait context.sync();\n\n// \"Maken profiel\" keeps its trailing \", \" un-struck, so target just the\n// \"Maken profiel\" span instead of the whole paragraph.\nconst makenProfiel = body.search(\"Maken profiel\", { matchCase: false }).getFirst();\nmakenProfiel.font.strikeThrough = true;\nawait context.sync();\n\n// Move the _GoBack bookmark off the old Facebook split and onto this range.\ncontext.document.deleteBookmark(\"_GoBack\");\nmakenProfiel.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// \"Faceb\" / \"ook\" were two separate runs (with the bookmark between them);\n// merge them into a single \"Facebook\" run.\nconst facebPart = body.search(\"Faceb\", { matchCase: true }).getFirst();\nconst ookPart = body.search(\"ook\", { matchCase: true }).getFirst();\nawait context.sync();\n\nconst facebookRange = facebPart.expandTo(ookPart);\nfacebookRange.insertText(\"Facebook\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Strike through the \"casual\" (Member/Visitor) use-case lines that are no\n# longer applicable, and tidy up the \"Facebook\" spell-checked word (which\n# used to be split into two runs with the _GoBack bookmark wedged between\n# them) into a single run, moving that bookmark onto \"Maken profiel\" instead.\n\n$d = $word.ActiveDocument\n\n# Paragraphs whose entire text gets struck through.\n$fullyStruckThrough = @(\n    \"Betalen contributie,  Fully-Dressed\",\n    \"Huren golfveld, \",\n    \"Opvragen statistieken,  Fully-Dressed \",\n    \"Opvragen Evenementen\",\n    \"Kopen artikelen,  Fully-Dressed \"\n)\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd(\"`r\", \"`n\")\n    if ($fullyStruckThrough -contains $text) {\n        $p.Range.Font.StrikeThrough = 1\n    }\n}\n\n# \"Maken profiel\" keeps its trailing \", \" un-struck, so target just the\n# \"Maken profiel\" span instead of the whole paragraph.\n$makenProfiel = $d.Content\n$makenProfiel.Find.Execute(\"Maken profiel\") | Out-Null\n$makenProfiel.Font.StrikeThrough = 1\n\n# Move the _GoBack bookmark off the old Facebook split and onto this range.\n$d.Bookmarks(\"_GoBack\").Delete()\n$makenProfiel2 = $d.Content\n$makenProfiel2.Find.Execute(\"Maken profiel\") | Out-Null\n$makenProfiel2.Bookmarks.Add(\"_GoBack\")\n\n# \"Faceb\" / \"ook\" were two separate runs (with the bookmark between them);\n# merge them into a single \"Facebook\" run.\n$facebPart = $d.Content\n$facebPart.Find.Execute(\"Faceb\") | Out-Null\n$start = $facebPart.Start\n\n$ookPart = $d.Content\n$ookPart.Find.Execute(\"ook\") | Out-Null\n$end = $ookPart.End\n\n$facebookRange = $d.Range($start, $end)\n$facebookRange.Text = \"Facebook\"\n"}
